$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '66.809.22'
$ws.Cells.Item(2, 5).Value = '  -4.21%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '3.382.66'
$ws.Cells.Item(3, 5).Value = '  -4.80%  '
$ws.Cells.Item(4, 5).Value = '  +0.26%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '563.18'
$ws.Cells.Item(5, 5).Value = '  -4.35%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '184.87'
$ws.Cells.Item(6, 5).Value = '  -7.05%  '
$ws.Cells.Item(7, 5).Value = '  -2.33%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  -0.01%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '3.375.03'
$ws.Cells.Item(9, 5).Value = '  -4.61%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.189'
$ws.Cells.Item(10, 5).Value = '  -8.83%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.599'
$ws.Cells.Item(11, 5).Value = '  -4.97%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '48.41'
$ws.Cells.Item(12, 5).Value = '  -7.52%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000272'
$ws.Cells.Item(13, 5).Value = '  -6.36%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '8.79'
$ws.Cells.Item(14, 5).Value = '  -6.12%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '3.929.59'
$ws.Cells.Item(15, 5).Value = '  -4.24%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '609.27'
$ws.Cells.Item(16, 5).Value = '  -11.18%  '
$ws.Cells.Item(17, 2).Value = 'Chainlink'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '18.32'
$ws.Cells.Item(17, 5).Value = '  -1.95%  '
$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '66.771.25'
$ws.Cells.Item(18, 5).Value = '  -4.24%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '3.383.61'
$ws.Cells.Item(19, 5).Value = '  -4.66%  '
$ws.Cells.Item(20, 5).Value = '  -3.01%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '11.68'
$ws.Cells.Item(21, 5).Value = '  -6.53%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '0.921'
$ws.Cells.Item(22, 5).Value = '  -5.46%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '17.07'
$ws.Cells.Item(23, 5).Value = '  -4.79%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '5.12'
$ws.Cells.Item(24, 5).Value = '  -2.41%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '98.93'
$ws.Cells.Item(25, 5).Value = '  -8.73%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '4.11'
$ws.Cells.Item(26, 5).Value = '  -7.02%  '
$ws.Cells.Item(27, 5).Value = '  +0.40%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.74'
$ws.Cells.Item(28, 5).Value = '  -7.46%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '9.55'
$ws.Cells.Item(29, 5).Value = '  -7.67%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '8.88'
$ws.Cells.Item(30, 5).Value = '  -8.77%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '30.97'
$ws.Cells.Item(31, 5).Value = '  -8.28%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.94'
$ws.Cells.Item(32, 5).Value = '  -10.59%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '6.38'
$ws.Cells.Item(33, 5).Value = '  -8.28%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '11.24'
$ws.Cells.Item(34, 5).Value = '  -6.05%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '552.93'
$ws.Cells.Item(35, 5).Value = '  +9.77%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '3.890.89'
$ws.Cells.Item(36, 5).Value = '  +2.22%  '
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.106'
$ws.Cells.Item(37, 5).Value = '  -5.39%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '58.51'
$ws.Cells.Item(38, 5).Value = '  -6.21%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.999'
$ws.Cells.Item(39, 5).Value = '  -0.04%  '
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '3.48'
$ws.Cells.Item(40, 5).Value = '  -6.68%  '
$ws.Cells.Item(41, 2).Value = 'CoreDAO'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.51'
$ws.Cells.Item(41, 5).Value = '  +29.44%  '
$ws.Cells.Item(42, 2).Value = 'PEPE'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.0₃0726'
$ws.Cells.Item(42, 5).Value = '  -11.59%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '2.73'
$ws.Cells.Item(43, 5).Value = '  -8.21%  '
$ws.Cells.Item(44, 5).Value = '  -5.56%  '
$ws.Cells.Item(45, 5).Value = '  -5.72%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '32.49'
$ws.Cells.Item(46, 5).Value = '  -7.27%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.0421'
$ws.Cells.Item(47, 5).Value = '  -8.77%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.28'
$ws.Cells.Item(48, 5).Value = '  -2.88%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '2.73'
$ws.Cells.Item(49, 5).Value = '  -8.05%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '0.131'
$ws.Cells.Item(50, 5).Value = '  -4.79%  '
$ws.Cells.Item(51, 5).Value = '  -0.01%  '
